$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("P16").Value = "W"
